# Commit: "add validation for ef config at end and rerun with new output"
#
# This inserts a new "CO2AI Taxonomy" column (H), shifting the former
# H:L columns (GBU, Category, SAP ID, Data Source, Very Custom) one to
# the right (I:M). The former "Taxonomy" column (F) is repurposed as a
# numeric "EF ID CO2" column, and a handful of category labels are
# re-cased / replaced with fuller taxonomy strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column before the old "GBU" column (H).
#    This shifts H->I, I->J, J->K, K->L, L->M automatically, carrying
#    all existing values/styles along with it.
$ws.Columns("H").Insert()

# 2. Re-purpose column F: rename header, and replace the old
#    "level 3" taxonomy strings with the new numeric EF ids
#    (only rows 2 and 3 keep a value; the rest become blank).
$ws.Range("F1").Value = "EF ID CO2"
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("F11").Value = ""

# 3. Populate the newly inserted column H ("CO2AI Taxonomy").
$ws.Range("H1").Value = "CO2AI Taxonomy"
$ws.Range("H8").Value = "Transportation `$ Goods transportation & distribution `$ Truck"
$ws.Range("H9").Value = "Transportation `$ Goods transportation & distribution `$ Truck"
$ws.Range("H10").Value = "Use and end of life (waste) `$ Waste treatment"
$ws.Range("H11").Value = "Use and end of life (waste) `$ Waste treatment"

# 4. Column J (old "Category" content, now shifted from old I) gets
#    re-cased / renamed category labels.
$ws.Range("J2").Value = "Glass sourcing"
$ws.Range("J3").Value = "Glass sourcing"
$ws.Range("J4").Value = "Ingredient sourcing"
$ws.Range("J5").Value = "Ingredient sourcing"
$ws.Range("J6").Value = "Ingredient sourcing"
$ws.Range("J7").Value = "Ingredient sourcing"
$ws.Range("J10").Value = "Use and end of life"
$ws.Range("J11").Value = "Use and end of life"
